$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.889831900596619
$ws.Range("B1").Value = 1.883208155632019
$ws.Range("C1").Value = 7.88377046585083
$ws.Range("D1").Value = 0.9500036835670471
$ws.Range("E1").Value = 0.4147926270961761
